$d = $word.ActiveDocument

# 1. Update the activation date.
$d.Content.Find.Execute(
    "Ativação: 01/01/2012", $true, $false, $false, $false, $false,
    $true, 1, $false, "Ativação: 01/01/2023", 2)

# 2. Insert an italic English translation paragraph right after the
#    Portuguese "Objetivos" paragraph.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Apresentar os conceitos de spintrônica e as potenciais aplicações em computação quântica.") {
        $target = $p
        break
    }
}
$target.Range.InsertParagraphAfter()
$newIndex = $target.Index + 1
$newPara = $d.Paragraphs.Item($newIndex)
$newRange = $newPara.Range
$text1 = "To present the concepts of spintronics and the potential applications in quantum computing."
$newRange.InsertAfter($text1)
$s1 = $newRange.Start
$e1 = $s1 + $text1.Length
$fmtRange1 = $d.Range($s1, $e1)
$fmtRange1.Font.Italic = $true

# 3. Insert an italic English translation paragraph right after the
#    "Programa resumido" body paragraph (the one with no <w:br/> breaks).
$target2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Introdução à nanotecnologia. Spintrônica de metais. Spintrônica de semicondutores. Dispositivos da spintrônica. Introdução à computação clássica. Introdução à computação quântica. Algoritmos quânticos. Decoerência. Pontos quânticos. Transistor de Kane. Introdução a modelos da consciência: o cérebro é um computador quântico?") {
        $target2 = $p
        break
    }
}
$target2.Range.InsertParagraphAfter()
$newIndex2 = $target2.Index + 1
$newPara2 = $d.Paragraphs.Item($newIndex2)
$newRange2 = $newPara2.Range
$text2 = "Introduction to nanotechnology. Metal spintronics. Semiconductor spintronics. Spintronics devices. Introduction to classical computing. Introduction to quantum computing. Quantum Algorithms. decoherence. Quantum Dots. Kane transistor. Introduction to models of consciousness: is the brain a quantum computer?"
$newRange2.InsertAfter($text2)
$s2 = $newRange2.Start
$e2 = $s2 + $text2.Length
$fmtRange2 = $d.Range($s2, $e2)
$fmtRange2.Font.Italic = $true

# 4. Collapse the "Programa" body paragraph's multiple <w:br/>-separated
#    runs into a single run (remove line breaks, keep text verbatim),
#    then add an italic English translation paragraph after it.
$oldPt = "Introdução à nanotecnologia." + [char]11 + `
    "Spintrônica de metais. Spintrônica de semicondutores " + [char]11 + `
    "Dispositivos da spintrônica." + [char]11 + `
    "Introdução à computação clássica. Introdução à computação quântica. Algoritmos quânticos." + [char]11 + `
    "Decoerência. Pontos quânticos. " + [char]11 + `
    "Transistor de Kane." + [char]11 + `
    "Introdução a modelos da consciência: o cérebro é um computador quântico?"
$mergedPt = "Introdução à nanotecnologia.Spintrônica de metais. Spintrônica de semicondutores Dispositivos da spintrônica.Introdução à computação clássica. Introdução à computação quântica. Algoritmos quânticos.Decoerência. Pontos quânticos. Transistor de Kane.Introdução a modelos da consciência: o cérebro é um computador quântico?"
$d.Content.Find.Execute($oldPt, $true, $false, $false, $false, $false, $true, 1, $false, $mergedPt, 2)

$target3 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $mergedPt) {
        $target3 = $p
        break
    }
}
$target3.Range.InsertParagraphAfter()
$newIndex3 = $target3.Index + 1
$newPara3 = $d.Paragraphs.Item($newIndex3)
$newRange3 = $newPara3.Range
$text3 = "Introduction to nanotechnology.Metal spintronics. Semiconductor SpintronicsSpintronics devices.Introduction to classical computing. Introduction to quantum computing. Quantum Algorithms.decoherence. Quantum Dots.Kane transistor.Introduction to models of consciousness: is the brain a quantum computer?"
$newRange3.InsertAfter($text3)
$s3 = $newRange3.Start
$e3 = $s3 + $text3.Length
$fmtRange3 = $d.Range($s3, $e3)
$fmtRange3.Font.Italic = $true
